$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value2 = 77650

# Row 3
$ws.Range("A3").Value2 = 112421694
$ws.Range("B3").Value2 = 78713
$ws.Range("D3").Value2 = "NT"
$ws.Range("E3").Value2 = 6458
$ws.Range("F3").Value2 = "Lunglav"
$ws.Range("G3").Value2 = "Lobaria pulmonaria"
$ws.Range("H3").Value2 = "(L.) Hoffm."
$ws.Range("I3").Value2 = ""
$ws.Range("Q3").Value2 = 491106
$ws.Range("R3").Value2 = 6954854
$ws.Range("Z3").Value2 = "15:54"
$ws.Range("AB3").Value2 = "15:54"
$ws.Range("AC3").Value2 = "Många tussar av Lunglav på en lång fallen sälg"

# Row 4
$ws.Range("A4").Value2 = 112422037
$ws.Range("B4").Value2 = 78713
$ws.Range("E4").Value2 = 6458
$ws.Range("F4").Value2 = "Lunglav"
$ws.Range("G4").Value2 = "Lobaria pulmonaria"
$ws.Range("H4").Value2 = "(L.) Hoffm."
$ws.Range("P4").Value2 = "Högberget (Högberget), Jmt"
$ws.Range("Q4").Value2 = 491071
$ws.Range("R4").Value2 = 6954842
$ws.Range("S4").Value2 = 1
$ws.Range("Z4").Value2 = "16:09"
$ws.Range("AB4").Value2 = "16:09"
$ws.Range("AC4").Value2 = "Annars nästan bara gran"

# Row 5
$ws.Range("A5").Value2 = 112421322
$ws.Range("B5").Value2 = 96735
$ws.Range("D5").Value2 = "VU"
$ws.Range("E5").Value2 = 220787
$ws.Range("F5").Value2 = "Knärot"
$ws.Range("G5").Value2 = "Goodyera repens"
$ws.Range("H5").Value2 = "(L.) R. Br."
$ws.Range("I5").Value2 = "'40"
$ws.Range("Q5").Value2 = 491128
$ws.Range("R5").Value2 = 6954848
$ws.Range("Z5").Value2 = "15:41"
$ws.Range("AB5").Value2 = "15:41"

# Row 6
$ws.Range("A6").Value2 = 112430267
$ws.Range("B6").Value2 = 77403
$ws.Range("E6").Value2 = 228912
$ws.Range("F6").Value2 = "Mörk kolflarnlav"
$ws.Range("G6").Value2 = "Carbonicola myrmecina"
$ws.Range("H6").Value2 = "(Ach.) Bendiksby & Timdal"
$ws.Range("P6").Value2 = "Högberget, Jmt"
$ws.Range("Q6").Value2 = 491026
$ws.Range("R6").Value2 = 6954834
$ws.Range("S6").Value2 = 5
$ws.Range("Z6").Value2 = ""
$ws.Range("AB6").Value2 = ""
$ws.Range("AC6").Value2 = "Kan också vara ickemörk kolflarnlav. Fanns på kolad mycket gammal stubbe."

# Row 7
$ws.Range("A7").Value2 = 112423544
$ws.Range("B7").Value2 = 78713
$ws.Range("Q7").Value2 = 491079
$ws.Range("R7").Value2 = 6954795
$ws.Range("Z7").Value2 = "16:55"
$ws.Range("AB7").Value2 = "16:55"
$ws.Range("AC7").Value2 = ""
